$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new "broker" column at the front (A), shifting B:F -> C:G ---
$ws.Columns("A:A").Insert()

# Copy the formatting from the (now-shifted) date column B onto the new
# column A so the header gets the bold header style (s=1) and the data rows
# get the same date-format style (s=2) the rest of the row already uses.
$ws.Range("B1:B20").Copy() | Out-Null
$ws.Range("A1:A20").PasteSpecial(-4122) | Out-Null

# --- Fill in the broker values for the existing 19 transactions ---
$ws.Range("A1").Value2 = "broker"
$ws.Range("A2:A20").Value2 = "company_a"

# --- Pre-format the two new rows by copying row 2's per-column formats down ---
$ws.Range("A2:G2").Copy() | Out-Null
$ws.Range("A21:G21").PasteSpecial(-4122) | Out-Null
$ws.Range("A22:G22").PasteSpecial(-4122) | Out-Null

# --- Append the two new company_b / NFLX transactions ---
$ws.Range("A21").Value2 = "company_b"
$ws.Range("B21").Value2 = 44701
$ws.Range("C21").Value2 = "Buy"
$ws.Range("D21").Value2 = "NFLX"
$ws.Range("E21").Value2 = 5
$ws.Range("F21").Value2 = 186.35

$ws.Range("A22").Value2 = "company_b"
$ws.Range("B22").Value2 = 44722
$ws.Range("C22").Value2 = "Sell"
$ws.Range("D22").Value2 = "NFLX"
$ws.Range("E22").Value2 = 5
$ws.Range("F22").Value2 = 182.94

# Assigning one formula to the whole G21:G22 block in one shot lets the
# engine register it as a single shared formula (mirrors Excel's fill-down
# behaviour), same as all the other cost-column formulas in this sheet.
$ws.Range("G21:G22").Formula = "=E21*F21"

# --- Cosmetic: selection moves to F23, matching the saved workbook state ---
$ws.Range("F23").Select() | Out-Null
